# Replace the Harris and d'Abramo (2015) flux data (column L, rows 2-45)
# with the corrected values, rerun/reformat the column, and update the
# active selection left over from working on the plot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = "2.2400000000000001E-9"
    3  = "2.2400000000000001E-9"
    4  = "2.2400000000000001E-9"
    5  = "2.2400000000000001E-9"
    6  = "4.4800000000000002E-9"
    7  = "4.4800000000000002E-9"
    8  = "6.7299999999999997E-9"
    9  = "1.1199999999999999E-8"
    10 = "1.7900000000000001E-8"
    11 = "3.3699999999999997E-8"
    12 = "5.84E-8"
    13 = "1.24E-7"
    14 = "2.4900000000000002E-7"
    15 = "4.4200000000000001E-7"
    16 = "7.3E-7"
    17 = "1.1400000000000001E-6"
    18 = "1.7799999999999999E-6"
    19 = "2.7300000000000001E-6"
    20 = "4.3100000000000002E-6"
    21 = "6.3500000000000002E-6"
    22 = "9.4399999999999994E-6"
    23 = "1.34E-5"
    24 = "1.91E-5"
    25 = "2.73E-5"
    26 = "3.8600000000000003E-5"
    27 = "5.66E-5"
    28 = "8.81E-5"
    29 = "1.4999999999999999E-4"
    30 = "2.8200000000000002E-4"
    31 = "5.9699999999999998E-4"
    32 = "1.4E-3"
    33 = "3.29E-3"
    34 = "7.6400000000000001E-3"
    35 = "1.6E-2"
    36 = "3.5700000000000003E-2"
    37 = "6.9699999999999998E-2"
    38 = "0.1303"
    39 = "0.23089999999999999"
    40 = "0.39739999999999998"
    41 = "0.57720000000000005"
    42 = "1.0345"
    43 = "1.6094999999999999"
    44 = "3.0413999999999999"
    45 = "8.1696000000000009"
}

# Rows 2-34 keep the scientific-notation display used for the rest of the
# exponential tail of the distribution; rows 35-45 are large enough that
# they were (re)pasted in with plain/general formatting instead.
foreach ($r in 2..34) {
    $cell = $ws.Range("L$r")
    $cell.NumberFormat = "0.00E+00"
    $cell.Value = [double]$newValues[$r]
    $cell.Font.Color = 0
}

foreach ($r in 35..45) {
    $cell = $ws.Range("L$r")
    $cell.Value = [double]$newValues[$r]
    $cell.Font.Color = 0
}

# Leave the selection where the author last clicked while reviewing the
# refreshed impactor flux plot.
$ws.Range("M26").Select() | Out-Null
